$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$aboutWs = $wb.Worksheets.Item("About")
$boundariesWs = $wb.Worksheets.Item("Boundaries and methane sources")

# Update the "About" sheet: A2 (Version line) and A6 (Recommended Citation)
$a2range = $aboutWs.Range("A2")
$a2text = $a2range.Value()
$a2range.Value = $a2text.Replace($oldStamp, $newStamp)

$a6range = $aboutWs.Range("A6")
$a6text = $a6range.Value()
$a6range.Value = $a6text.Replace($oldStamp, $newStamp)

# Update the "Boundaries and methane sources" sheet: build_version column (S), rows 2-24
$usedRange = $boundariesWs.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $boundariesWs.Cells.Item($r, 19)  # column S = 19
    $cellValue = $cell.Value()
    if ($cellValue -ne $null) {
        $cellText = $cellValue.ToString()
        if ($cellText.Contains($oldStamp)) {
            $cell.Value = $cellText.Replace($oldStamp, $newStamp)
        }
    }
}
